$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (rows 2-31): Timestamp (A), IGCC Import MW (B), IGCC Export MW (C)
$data = @(
    @(46006, 0, 11.364),
    @(46006.01041666666, 0.005, 11.286),
    @(46006.02083333334, 0, 11.7),
    @(46006.03125, 0, 20.788),
    @(46006.04166666666, 0, 23.581),
    @(46006.05208333334, 0, 10.778),
    @(46006.0625, 0, 12.671),
    @(46006.07291666666, 2.307, 1.115),
    @(46006.08333333334, 3.471, 5.077),
    @(46006.09375, 0.367, 6.033),
    @(46006.10416666666, 0.032, 6.226),
    @(46006.11458333334, 2.972, 0.499),
    @(46006.125, 0.5679999999999999, 3.026),
    @(46006.13541666666, 0.103, 7.661),
    @(46006.14583333334, 0.124, 4.172),
    @(46006.15625, 0, 2.709),
    @(46006.16666666666, 0.075, 15.471),
    @(46006.17708333334, 0.992, 6.021),
    @(46006.1875, 1.401, 1.025),
    @(46006.19791666666, 0.718, 3.187),
    @(46006.20833333334, 0, 21.889),
    @(46006.21875, 0.289, 8.550000000000001),
    @(46006.22916666666, 0.499, 1.314),
    @(46006.23958333334, 5.47, 1.389),
    @(46006.25, 0.117, 13.95),
    @(46006.26041666666, 0.234, 1.93),
    @(46006.27083333334, 0.238, 0.31),
    @(46006.28125, 0, 4.544),
    @(46006.29166666666, 0.444, 2.286),
    @(46006.30208333334, 0, 0)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the now-unused trailing rows (previously up to row 44)
$lastRow = 44
$deleteRange = $ws.Range("A32:C$lastRow")
$deleteRange.EntireRow.Delete()
